$wb = $excel.ActiveWorkbook

# Registration sheet: A2 email address changed from zqio@test.com to ypqh@test.com
$wsReg = $wb.Worksheets.Item("Registration")
$wsReg.Range("A2").Value = "ypqh@test.com"

# test_suite sheet: Runmode column (B) for Login, Parameter, VerifyLoginPage rows changed from N to Y
$wsSuite = $wb.Worksheets.Item("test_suite")
$wsSuite.Range("B2").Value = "Y"
$wsSuite.Range("B4").Value = "Y"
$wsSuite.Range("B5").Value = "Y"

# Update selection state on test_suite sheet
$wsSuite.Range("B2:B6").Select()
